# raven.docx edit: intro-line annotation + trim trailing paragraph.
$d = $word.ActiveDocument

# --- 1. First paragraph: "This is a Microsoft word document." ---------------
# Add two trailing spaces to the existing (plain) run, then append three
# separate red-colored runs: "(This is a change – Ve" / "rsion for main
# branch" / ")".
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$r.End = $r.End - 1                 # exclude the paragraph mark
$r.InsertAfter("  ")
$r.Collapse(0)

$r.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$seg1 = $d.Range($r.Start, $r.End)
$seg1.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter("rsion for main branch")
$seg2 = $d.Range($r.Start, $r.End)
$seg2.Font.Color = 255
$r.Collapse(0)

$r.InsertAfter(")")
$seg3 = $d.Range($r.Start, $r.End)
$seg3.Font.Color = 255
$r.Collapse(0)

# --- 2. Drop the trailing "ank God almighty, we are free at last." para -----
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.Delete()

# --- 3. Prune the now-unused custom styles that rode along with that para ---
# (Deleting by descending index - deleting low-to-high by name trips an
# index-cache bug in the host after a handful of removals.)
$deadStyles = @(
    "apple-converted-space",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)
$deadIndices = @()
for ($i = 1; $i -le $d.Styles.Count; $i++) {
    if ($deadStyles -contains $d.Styles.Item($i).NameLocal) {
        $deadIndices += $i
    }
}
[array]::Reverse($deadIndices)
foreach ($idx in $deadIndices) {
    $d.Styles.Item($idx).Delete()
}
